# Regenerate save_data: switch column G (header "K") from the old
# Strike# values to the recalculated K values (std/mean based s_vals).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..71 (one entry per data row, in order).
$kValues = @(
    2, 2, 0, 1, 1, 1, 0, 0, 2, 2,
    3, 2, 0, 0, 1, 2, 0, 2, 0, 0,
    1, 3, 0, 3, 0, 2, 2, 2, 0, 1,
    1, 1, 1, 3, 1, 0, 0, 1, 2, 1,
    1, 2, 0, 2, 1, 1, 0, 4, 2, 0,
    3, 1, 1, 0, 1, 0, 2, 1, 1, 2,
    1, 5, 2, 1, 4, 3, 3, 2, 3, 1
)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
